$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update "Version" value (row 3) and "Date" value (row 8)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before the current row 11 ("Description") for the new
# "Jurisdiction" property, shifting everything below down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"

Write-Output "done"
